$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A27").Value = "Record"
$ws.Range("B27").Value = "Balanço Geral"
$ws.Range("C27").Value = "Infraestrutura"
$ws.Range("D27").Value = "2025-04-01T13:16"
$ws.Range("E27").Value = "Negativo"
$ws.Range("F27").Value = "Bueiro sem tampa em rua do Pq. Tomás Coelho gera transtornos para motoristas. Bueiro sem tampa já teria procurado acidentes. *Com nota*. *nota coberta*"
